$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "./RESULTS/final scores/trial 3/"
$new = "./RESULTS/trial 3/trial 3 results/"

$cells = @("B2", "C2", "D2", "E2", "F2", "G2", "H2", "I2", "J2")

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $text = [string]$cell.Value()
    $text = $text.Replace($old, $new)
    $cell.Value = $text
}
